$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 160, shifting existing rows 160-230 down to 161-231
$ws.Rows.Item(160).Insert()

# Populate the new row 160 with its values
$ws.Cells.Item(160, 1).Value = 5
$ws.Cells.Item(160, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(160, 3).Value = "Maule"
$ws.Cells.Item(160, 4).Value = 44523
$ws.Cells.Item(160, 4).NumberFormat = $ws.Cells.Item(161, 4).NumberFormat
$ws.Cells.Item(160, 5).Value = 7
$ws.Cells.Item(160, 6).Value = 100112023
$ws.Cells.Item(160, 7).Value = "Brócoli"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 5000
$ws.Cells.Item(160, 11).Value = 500
$ws.Cells.Item(160, 12).Value = 500
$ws.Cells.Item(160, 13).Value = 500
$ws.Cells.Item(160, 14).Value = "$/unidad"
$ws.Cells.Item(160, 15).Value = "Región del Maule"
$ws.Cells.Item(160, 16).Value = 500
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"
